# Re-orders the monthly data rows so that, within each contiguous
# "year block" in column A, the Oct/Nov/Dec rows move to the front of
# the block (carrying their own A:E values with them), while the other
# months keep their original relative order. Also drops column F
# (the "...出口交货值" absolute-value column) entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow  = 53
$numRows      = $lastDataRow - $firstDataRow + 1
$numCols      = 5

# --- 1. Snapshot every data row (A:E) before touching anything -------------
# Range.Value2 on a multi-cell range comes back as a 1-based 2D SAFEARRAY:
# indices run [1..numRows, 1..numCols].
$snapshot = $ws.Range("A" + $firstDataRow + ":E" + $lastDataRow).Value2

$years  = @()
$months = @()
for ($i = 1; $i -le $numRows; $i++) {
    $label = $snapshot[$i, 1]
    $years  += $label.Substring(0, 4)
    $months += $label.Substring(5, 2)
}

# --- 2. Split the row indices (1-based within the snapshot) into
#        contiguous blocks that share the same "year" label -----------------
$blockStarts = @()
$blockEnds   = @()
$curStart = 1
for ($i = 2; $i -le $numRows; $i++) {
    if ($years[$i - 1] -ne $years[$i - 2]) {
        $blockStarts += $curStart
        $blockEnds   += ($i - 1)
        $curStart = $i
    }
}
$blockStarts += $curStart
$blockEnds   += $numRows

# --- 3. For every block, move the Oct/Nov/Dec rows to the front,
#        preserving the relative order within both groups -------------------
$newOrder = @()
for ($b = 0; $b -lt $blockStarts.Count; $b++) {
    $s = $blockStarts[$b]
    $e = $blockEnds[$b]

    $tail = @()
    $head = @()
    for ($i = $s; $i -le $e; $i++) {
        $m = $months[$i - 1]
        if ($m -eq "10" -or $m -eq "11" -or $m -eq "12") {
            $tail += $i
        } else {
            $head += $i
        }
    }
    $newOrder += $tail
    $newOrder += $head
}

# --- 4. Build the rearranged buffer (0-based .NET array, as required by
#        `New-Object 'object[,]'`) and write it back in one shot ------------
$buf = New-Object 'object[,]' $numRows, $numCols
for ($i = 0; $i -lt $numRows; $i++) {
    $src = $newOrder[$i]
    for ($j = 0; $j -lt $numCols; $j++) {
        $buf[$i, $j] = $snapshot[$src, $j + 1]
    }
}
$ws.Range("A" + $firstDataRow + ":E" + $lastDataRow).Value = $buf

# --- 5. Drop column F entirely (header + data) ------------------------------
$ws.Columns.Item(6).Delete()
